$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2201492537313433
$ws.Range("C2").Value = 0.4925373134328358
$ws.Range("J2").Value = 0.01865671641791045
$ws.Range("P2").Value = 0.1529850746268657
$ws.Range("S2").Value = 0.1156716417910448
$ws.Range("B3").Value = 0.006711409395973154
$ws.Range("C3").Value = 0.02013422818791946
$ws.Range("J3").Value = 0.05369127516778523
$ws.Range("P3").Value = 0.7986577181208053
$ws.Range("S3").Value = 0.1208053691275168
$ws.Range("J4").Value = 0.075
$ws.Range("P4").Value = 0.575
$ws.Range("S4").Value = 0.35
$ws.Range("B6").Value = 0.0582010582010582
$ws.Range("D6").Value = 0.01058201058201058
$ws.Range("F6").Value = 0.03703703703703703
$ws.Range("J6").Value = 0.3227513227513227
$ws.Range("O6").Value = 0.02645502645502645
$ws.Range("Q6").Value = 0.1428571428571428
$ws.Range("R6").Value = 0.07936507936507936
$ws.Range("S6").Value = 0.3227513227513227
$ws.Range("B7").Value = 0.1304347826086956
$ws.Range("D7").Value = 0.02484472049689441
$ws.Range("F7").Value = 0.04968944099378882
$ws.Range("J7").Value = 0.1180124223602484
$ws.Range("O7").Value = 0.006211180124223602
$ws.Range("Q7").Value = 0.1118012422360248
$ws.Range("R7").Value = 0.06832298136645963
$ws.Range("S7").Value = 0.4906832298136646
$ws.Range("B8").Value = 0.09414758269720101
$ws.Range("D8").Value = 0.01526717557251908
$ws.Range("F8").Value = 0.06615776081424936
$ws.Range("J8").Value = 0.1170483460559796
$ws.Range("O8").Value = 0.01526717557251908
$ws.Range("Q8").Value = 0.1679389312977099
$ws.Range("R8").Value = 0.08905852417302799
$ws.Range("S8").Value = 0.4351145038167939
$ws.Range("B9").Value = 0.08035714285714286
$ws.Range("D9").Value = 0.008928571428571428
$ws.Range("F9").Value = 0.1071428571428571
$ws.Range("J9").Value = 0.1607142857142857
$ws.Range("O9").Value = 0.03571428571428571
$ws.Range("Q9").Value = 0.125
$ws.Range("R9").Value = 0.09821428571428571
$ws.Range("S9").Value = 0.3839285714285715
$ws.Range("B10").Value = 0.1155638397017707
$ws.Range("D10").Value = 0.02423112767940354
$ws.Range("E10").Value = 0.002795899347623486
$ws.Range("F10").Value = 0.07269338303821063
$ws.Range("J10").Value = 0.1314072693383038
$ws.Range("O10").Value = 0.0260950605778192
$ws.Range("Q10").Value = 0.2087604846225536
$ws.Range("R10").Value = 0.08480894687791239
$ws.Range("S10").Value = 0.3336439888164026
$ws.Range("F11").Value = 0.008888888888888889
$ws.Range("G11").Value = 0.1511111111111111
$ws.Range("J11").Value = 0.08
$ws.Range("K11").Value = 0.1911111111111111
$ws.Range("L11").Value = 0.5555555555555556
$ws.Range("S11").Value = 0.01333333333333333
$ws.Range("G12").Value = 0.7674418604651163
$ws.Range("J12").Value = 0.1705426356589147
$ws.Range("K12").Value = 0.007751937984496124
$ws.Range("L12").Value = 0.0310077519379845
$ws.Range("S12").Value = 0.02325581395348837
$ws.Range("G13").Value = 0.7169811320754716
$ws.Range("J13").Value = 0.2075471698113208
$ws.Range("S13").Value = 0.07547169811320754
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.02298850574712644
$ws.Range("H15").Value = 0.1436781609195402
$ws.Range("I15").Value = 0.05747126436781609
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("K15").Value = 0.05747126436781609
$ws.Range("M15").Value = 0.01149425287356322
$ws.Range("O15").Value = 0.04022988505747126
$ws.Range("S15").Value = 0.3333333333333333
$ws.Range("F16").Value = 0.02857142857142857
$ws.Range("H16").Value = 0.1542857142857143
$ws.Range("I16").Value = 0.04571428571428571
$ws.Range("J16").Value = 0.4171428571428571
$ws.Range("K16").Value = 0.1657142857142857
$ws.Range("M16").Value = 0.02285714285714286
$ws.Range("O16").Value = 0.05142857142857143
$ws.Range("S16").Value = 0.1142857142857143
$ws.Range("F17").Value = 0.01404494382022472
$ws.Range("H17").Value = 0.1685393258426966
$ws.Range("I17").Value = 0.06179775280898876
$ws.Range("J17").Value = 0.4353932584269663
$ws.Range("K17").Value = 0.07865168539325842
$ws.Range("M17").Value = 0.03370786516853932
$ws.Range("O17").Value = 0.0702247191011236
$ws.Range("S17").Value = 0.1376404494382023
$ws.Range("F18").Value = 0.005952380952380952
$ws.Range("H18").Value = 0.1369047619047619
$ws.Range("I18").Value = 0.05952380952380952
$ws.Range("J18").Value = 0.4880952380952381
$ws.Range("K18").Value = 0.09523809523809523
$ws.Range("M18").Value = 0.02380952380952381
$ws.Range("N18").Value = 0.005952380952380952
$ws.Range("O18").Value = 0.04166666666666666
$ws.Range("S18").Value = 0.1428571428571428
$ws.Range("F19").Value = 0.01459143968871595
$ws.Range("H19").Value = 0.2529182879377432
$ws.Range("I19").Value = 0.061284046692607
$ws.Range("J19").Value = 0.3608949416342412
$ws.Range("K19").Value = 0.09922178988326848
$ws.Range("M19").Value = 0.0301556420233463
$ws.Range("N19").Value = 0.0009727626459143969
$ws.Range("O19").Value = 0.06906614785992218
$ws.Range("S19").Value = 0.1108949416342412
